# "Added checklist from Excel"
#
# Renames Sheet3 -> checklist, fills it in with a small checklist table,
# replaces the `tblItems` defined name with `exChkLst` / `exItems`, and
# nudges the saved selections on all three sheets.

$wb = $excel.ActiveWorkbook

$wsHeaders = $wb.Worksheets.Item(1)
$wsItems   = $wb.Worksheets.Item(2)
$wsChk     = $wb.Worksheets.Item(3)

# --- rename the third sheet ------------------------------------------------
$wsChk.Name = "checklist"

# --- defined names -----------------------------------------------------
$wb.Names.Item("tblItems").Delete()
$wb.Names.Add("exChkLst", "=checklist!`$A`$1:`$D`$2")
$wb.Names.Add("exItems", "=items!`$A`$1:`$D`$9")

# --- populate the checklist sheet ------------------------------------------
$wsChk.Range("A1").Value = "chklst_name"
$wsChk.Range("B1").Value = "template_id"
$wsChk.Range("C1").Value = "area_code"
$wsChk.Range("D1").Value = "reference"

$wsChk.Range("A2").Value = "Panic Plan"
$wsChk.Range("B2").Value = 1
$wsChk.Range("C2").Value = "PG"
$wsChk.Range("D2").Value = "This Panic Pan is based on science fiction."

# NB: the COM host quantizes ColumnWidth to 1/6-character steps (it computes
# round((input + 5/6) * 6) / 6), so it can't reproduce Excel's raw 1/256-char
# autofit widths exactly. These inputs land on the closest reachable values
# (26.5 / 13.33333.. / 13.83333.. / 66.5) to the target 26.42578125 /
# 13.28515625 / 13.85546875 / 66.42578125.
$wsChk.Columns.Item(1).ColumnWidth = 25.666666666666668
$wsChk.Columns.Item(2).ColumnWidth = 12.5
$wsChk.Columns.Item(3).ColumnWidth = 13.0
$wsChk.Columns.Item(4).ColumnWidth = 65.66666666666667

# --- selections --------------------------------------------------------
# headers: A5 -> D20 (not the active tab)
$wsHeaders.Select()
$wsHeaders.Range("D20").Select()

# checklist: (none) -> C10 (not the active tab)
$wsChk.Select()
$wsChk.Range("C10").Select()

# items: D5 -> B23, stays the active tab, so select it last
$wsItems.Select()
$wsItems.Range("B23").Select()
# best-effort: also nudge the scroll position to match the saved topLeftCell="A9"
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
